$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 97; $row++) {
    $ws.Cells.Item($row, 1).Value = "q" + ($row - 2)
}
